$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Results for the "380 kV" case: bus voltage magnitudes (vm_pu) for rows 2-25.
# Columns B-F, I-N carry data (G stays fixed at 1, H/A are untouched).
$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

$data = @(
    @(1.02, 1.060575180405114, 1.059159839341727, 1.064270884879515, 1.068675108853807, 1.040471598234718, 1.065555321121076, 1.061890102576414, 1.066987273511069, 1.071379656788373, 1.067068531897971),
    @(1.02, 1.062570914073163, 1.060680131262395, 1.066219105108157, 1.070455424447177, 1.040928466747826, 1.067199763563828, 1.063222782717542, 1.068747830242324, 1.07297360315938, 1.068715309637612),
    @(1.02, 1.063858037923618, 1.061659937631226, 1.067475868775574, 1.071603648083809, 1.041220975459779, 1.068259378602209, 1.064080697198362, 1.069882756359328, 1.07400078825747, 1.069776429450888),
    @(1.02, 1.064398150960647, 1.062070925618544, 1.068003308846263, 1.072085480691203, 1.041343206471044, 1.068703795051201, 1.064440322483965, 1.070358877060199, 1.074431628445848, 1.070221477022212),
    @(1.02, 1.064488780716299, 1.062139878695019, 1.068091816015396, 1.07216633142947, 1.041363686403479, 1.068778353769565, 1.064500644510073, 1.070438761640625, 1.074503910974244, 1.070296141622508),
    @(1.02, 1.063865258813788, 1.061665432875193, 1.067482919975552, 1.071610089788512, 1.041222611614317, 1.068265320999163, 1.064085506601682, 1.069889122217857, 1.074006549024302, 1.069782380286726),
    @(1.02, 1.061250541981339, 1.05967445183397, 1.064930108597146, 1.069277565064015, 1.040626648024773, 1.066112001137761, 1.062341412601073, 1.067583158781484, 1.071919220078172, 1.067626002464133),
    @(1.02, 1.056609472365662, 1.056135272804645, 1.060401135392446, 1.065137685389502, 1.039552337097744, 1.062282622584394, 1.059233536083131, 1.0634861155625, 1.06820805288259, 1.063791185753983),
    @(1.02, 1.053491340635066, 1.053754028223525, 1.057359855429573, 1.062356577441896, 1.038819494107827, 1.059704965840331, 1.057137370370579, 1.060730834348847, 1.065710576666895, 1.061209868441683),
    @(1.02, 1.052135107891445, 1.052717513258101, 1.056037424910571, 1.061147023610321, 1.038498128988422, 1.058582666934049, 1.056223728839274, 1.05953180194891, 1.064623339755092, 1.060085975742224),
    @(1.02, 1.051630403049718, 1.0523316701694, 1.055545356439482, 1.060696918726071, 1.03837814494299, 1.058164846108527, 1.055883442320573, 1.059085504558988, 1.06421859558434, 1.059667561563206),
    @(1.02, 1.051738706941064, 1.052414472959214, 1.05565094610359, 1.060793505372922, 1.038403909873645, 1.058254513430569, 1.055956476857787, 1.059181279059021, 1.064305455493884, 1.05975735622313),
    @(1.02, 1.052093408168909, 1.052685636509693, 1.055996768056872, 1.061109834724461, 1.038488223656132, 1.058548149235078, 1.056195619519764, 1.059494929851734, 1.064589901885183, 1.060051409024162),
    @(1.02, 1.052311825919181, 1.052852598008412, 1.056209725745203, 1.061304625930419, 1.038540090464642, 1.058728941461393, 1.05634284078454, 1.059688057337462, 1.064765039287889, 1.060232457996175),
    @(1.02, 1.053581219168998, 1.053822702327823, 1.057447501838316, 1.062436737500034, 1.038840736280068, 1.059779317363166, 1.057197877911507, 1.060810282089034, 1.065782608600403, 1.061284325552209),
    @(1.02, 1.054375832941446, 1.05442975727482, 1.058222424838713, 1.063145441095582, 1.039028236507093, 1.060436525240611, 1.057732603115165, 1.061512605930298, 1.066419330728364, 1.061942466740222),
    @(1.02, 1.05483873487956, 1.054783320599357, 1.058673892001311, 1.063558304666727, 1.039137213024627, 1.060819271105819, 1.058043922878632, 1.061921684535945, 1.066790160623862, 1.06232575614846),
    @(1.02, 1.05499647449097, 1.054903788720072, 1.058827741162726, 1.063698994483183, 1.039174305423574, 1.060949677802285, 1.058149977698471, 1.062061072998354, 1.066916509838716, 1.062456348037393),
    @(1.02, 1.054290638897614, 1.054364680121804, 1.058139338189619, 1.063069456961043, 1.039008159819961, 1.060366074500275, 1.057675291923502, 1.061437312887053, 1.066351074488782, 1.061871915951754),
    @(1.02, 1.051988983670931, 1.052605808830107, 1.055894956040634, 1.061016706521608, 1.038463412382102, 1.058461707202751, 1.056125223505593, 1.059402593205556, 1.064506164428956, 1.059964844234234),
    @(1.02, 1.050536390686972, 1.051495093374085, 1.054478840871351, 1.059721288553989, 1.038117347411755, 1.057258852846637, 1.055145304756588, 1.058117929037968, 1.06334100099294, 1.058760281686965),
    @(1.02, 1.051306964194055, 1.052084370807965, 1.055230031814135, 1.060408474626897, 1.038301143168503, 1.057897038461093, 1.05566529004497, 1.058799469977344, 1.06395917625828, 1.059399373598194),
    @(1.02, 1.054329136218925, 1.054394087297064, 1.058176883113958, 1.063103792499812, 1.039017232813066, 1.06039790998701, 1.057701190164582, 1.061471336348356, 1.066381918273127, 1.061903796648532),
    @(1.02, 1.057813434721798, 1.057053996810798, 1.061575751469146, 1.066211586301647, 1.039832974643169, 1.063276877208773, 1.060041195490947, 1.064549418886695, 1.069171510779641, 1.064786852333971)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $cols.Count; $c++) {
        $ws.Cells.Item($r + 2, $cols[$c]).Value = $rowVals[$c]
    }
}